$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ((Intercept))
$ws.Range("B2").Value = 5.55954157615657
$ws.Range("C2").Value = 0.302997798995133
$ws.Range("D2").Value = 18.3484553174786
$ws.Range("E2").Value = 0.00000000000000000000000000000000126938171766596

# Row 3 (depression_mc)
$ws.Range("B3").Value = 0.0366148925537197
$ws.Range("C3").Value = 0.610755412255824
$ws.Range("D3").Value = 0.0599501728825991
$ws.Range("E3").Value = 0.952316172150894

# Row 4 (anhedonia_mc)
$ws.Range("B4").Value = -0.899849338922523
$ws.Range("C4").Value = 0.606242158905168
$ws.Range("D4").Value = -1.48430676703116
$ws.Range("E4").Value = 0.140904579927302

# Row 5 (depression_mc:anhedonia_mc)
$ws.Range("B5").Value = 0.690663905831003
$ws.Range("C5").Value = 1.21972185503999
$ws.Range("D5").Value = 0.566247052946638
$ws.Range("E5").Value = 0.57250723343244
